$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.042.58"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.12%  "
$ws.Range("D3").Value = "'3.124.38"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'591.92"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.81%  "
$ws.Range("D6").Value = "'136.31"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.13%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'3.109.36"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").Value = "'0.514"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.87%  "
$ws.Range("E10").Value = "  -3.95%  "
$ws.Range("D11").Value = "'5.28"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("D12").Value = "'0.454"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.98%  "
$ws.Range("D13").Value = "'0.0000245"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.57%  "
$ws.Range("D14").Value = "'33.96"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.00%  "
$ws.Range("D15").Value = "'3.637.05"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("D17").Value = "'63.049.00"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("D18").Value = "'3.120.29"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").Value = "'6.64"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.05%  "
$ws.Range("D20").Value = "'472.13"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("D21").Value = "'14.08"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.04%  "
$ws.Range("D22").Value = "'0.694"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.29%  "
$ws.Range("D23").Value = "'7.71"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("D24").Value = "'86.13"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").Value = "'12.86"
$ws.Range("D25").ClearFormats()
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'7.90"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.18%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").Value = "'6.92"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.41%  "
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").Value = "'26.55"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("D33").Value = "'0.107"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -7.34%  "
$ws.Range("E34").Value = "  -4.39%  "
$ws.Range("E35").Value = "  -2.54%  "
$ws.Range("D36").Value = "'5.76"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.54%  "
$ws.Range("D37").Value = "'52.09"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("D38").Value = "'0.0₃0696"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -8.42%  "
$ws.Range("D39").Value = "'421.22"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.52%  "
$ws.Range("D40").Value = "'0.0385"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("D41").Value = "'8.17"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'2.893.04"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.67"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -11.61%  "
$ws.Range("E44").Value = "  -6.05%  "
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'2.11"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.69%  "
$ws.Range("D48").Value = "'25.37"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.16%  "
$ws.Range("D50").Value = "'2.25"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.77%  "
$ws.Range("D51").Value = "'120.13"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.45%  "
